$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update the time_taken (column F) timestamps on the "data" sheet to reflect a later run
$ws.Range("F2").Value = "2021-10-05 14:34:54.729591"
$ws.Range("F3").Value = "2021-10-05 14:34:54.729598"
$ws.Range("F4").Value = "2021-10-05 14:34:54.729601"
$ws.Range("F5").Value = "2021-10-05 14:34:54.729604"
$ws.Range("F6").Value = "2021-10-05 14:34:54.729607"
$ws.Range("F7").Value = "2021-10-05 14:34:54.729609"
$ws.Range("F8").Value = "2021-10-05 14:34:54.729612"
$ws.Range("F9").Value = "2021-10-05 14:34:54.729614"
$ws.Range("F10").Value = "2021-10-05 14:34:54.729617"
$ws.Range("F11").Value = "2021-10-05 14:34:54.729620"
$ws.Range("F12").Value = "2021-10-05 14:34:54.729622"
$ws.Range("F13").Value = "2021-10-05 14:34:54.729624"
$ws.Range("F14").Value = "2021-10-05 14:34:54.729627"
$ws.Range("F15").Value = "2021-10-05 14:34:54.729629"
$ws.Range("F16").Value = "2021-10-05 14:34:54.729632"
$ws.Range("F17").Value = "2021-10-05 14:34:54.729635"
$ws.Range("F18").Value = "2021-10-05 14:34:54.729637"
$ws.Range("F19").Value = "2021-10-05 14:34:54.729640"
$ws.Range("F20").Value = "2021-10-05 14:34:54.729642"
$ws.Range("F21").Value = "2021-10-05 14:34:54.729645"
$ws.Range("F22").Value = "2021-10-05 14:34:54.729647"
$ws.Range("F23").Value = "2021-10-05 14:34:54.729650"
$ws.Range("F24").Value = "2021-10-05 14:34:54.729652"
$ws.Range("F25").Value = "2021-10-05 14:34:54.729655"
$ws.Range("F26").Value = "2021-10-05 14:34:54.729657"
$ws.Range("F27").Value = "2021-10-05 14:34:54.729660"
$ws.Range("F28").Value = "2021-10-05 14:34:54.729662"
$ws.Range("F29").Value = "2021-10-05 14:34:54.729665"
$ws.Range("F30").Value = "2021-10-05 14:34:54.729667"
$ws.Range("F31").Value = "2021-10-05 14:34:54.729670"
$ws.Range("F32").Value = "2021-10-05 14:34:54.729672"
$ws.Range("F33").Value = "2021-10-05 14:34:54.729675"
$ws.Range("F34").Value = "2021-10-05 14:34:54.729678"
$ws.Range("F35").Value = "2021-10-05 14:34:54.729680"
$ws.Range("F36").Value = "2021-10-05 14:34:54.729683"
$ws.Range("F37").Value = "2021-10-05 14:34:54.729685"
$ws.Range("F38").Value = "2021-10-05 14:34:54.729688"
$ws.Range("F39").Value = "2021-10-05 14:34:54.729690"
$ws.Range("F40").Value = "2021-10-05 14:34:54.729692"
$ws.Range("F41").Value = "2021-10-05 14:34:54.729695"
$ws.Range("F42").Value = "2021-10-05 14:34:54.729698"
$ws.Range("F43").Value = "2021-10-05 14:34:54.729700"
$ws.Range("F44").Value = "2021-10-05 14:34:54.729703"
$ws.Range("F45").Value = "2021-10-05 14:34:54.729705"
$ws.Range("F46").Value = "2021-10-05 14:34:54.729708"
$ws.Range("F47").Value = "2021-10-05 14:34:54.729710"
$ws.Range("F48").Value = "2021-10-05 14:34:54.729713"
$ws.Range("F49").Value = "2021-10-05 14:34:54.729715"
$ws.Range("F50").Value = "2021-10-05 14:34:54.729717"
$ws.Range("F51").Value = "2021-10-05 14:34:54.729720"
$ws.Range("F52").Value = "2021-10-05 14:34:54.729722"
$ws.Range("F53").Value = "2021-10-05 14:34:54.729725"
$ws.Range("F54").Value = "2021-10-05 14:34:54.729728"
$ws.Range("F55").Value = "2021-10-05 14:34:54.729730"
$ws.Range("F56").Value = "2021-10-05 14:34:54.729733"
$ws.Range("F57").Value = "2021-10-05 14:34:54.729735"
$ws.Range("F58").Value = "2021-10-05 14:34:54.729738"
$ws.Range("F59").Value = "2021-10-05 14:34:54.729740"
$ws.Range("F60").Value = "2021-10-05 14:34:54.729742"
$ws.Range("F61").Value = "2021-10-05 14:34:54.729745"
$ws.Range("F62").Value = "2021-10-05 14:34:54.729747"
$ws.Range("F63").Value = "2021-10-05 14:34:54.729750"
$ws.Range("F64").Value = "2021-10-05 14:34:54.729752"
$ws.Range("F65").Value = "2021-10-05 14:34:54.729755"
$ws.Range("F66").Value = "2021-10-05 14:34:54.729758"
$ws.Range("F67").Value = "2021-10-05 14:34:54.729761"
$ws.Range("F68").Value = "2021-10-05 14:34:54.729764"
$ws.Range("F69").Value = "2021-10-05 14:34:54.729766"
$ws.Range("F70").Value = "2021-10-05 14:34:54.729769"
$ws.Range("F71").Value = "2021-10-05 14:34:54.729771"
$ws.Range("F72").Value = "2021-10-05 14:34:54.729773"

# Add a new "metadata" worksheet after "data", describing the panel query itself
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "metadata"

# Copy the bold/bordered header style used on the "data" sheet for the header row and index column
$ws.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row describing the panel that was queried
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Monogenic Diabetes"
$newSheet.Range("C2").Value = 3093

# Keep "0.23" as text (matches source data), not a number
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.23"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").Value = "2021-08-22T04:36:52.626513Z"
$newSheet.Range("F2").Value = "2021-10-05 14:34:54.725973"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3093/?format=json"

# Keep "data" as the active/selected sheet, as it was before the edit
$ws.Activate()

Write-Output "metadata sheet added and timestamps refreshed"
